$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "257.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.07%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.33%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.518"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-6.04%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.611"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.69%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8505"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.31%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9276"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.71%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1378"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.08%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04217"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.25%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07001"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.41%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03050"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.62%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09099"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.76%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001539"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006028"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.21%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006019"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.34%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.469"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.40%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.174"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.49%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.212"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.33%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3080"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.77%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.92%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.903"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.24%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04261"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.19%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.52%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.29%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.11%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.98%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03798"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.19%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1100"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.03%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003893"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-37.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002389"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.03%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01402"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "32.32%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005348"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.79%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05099"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-42.42%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10,506.64%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
